$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197588562965393
$ws.Range("B1").Value = 2.484435796737671
$ws.Range("C1").Value = 4.242732048034668
$ws.Range("D1").Value = 2.08942174911499
$ws.Range("E1").Value = 1.182473301887512
